# Weekly update: prepend two new rows of data (new reporting week) above the
# existing data, shifting rows 15-45 down to 17-47, then populate the two new
# rows (15 and 16) with the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 15; existing row 15 (and below) shift down.
$ws.Rows.Item(15).Insert()
$ws.Rows.Item(15).Insert()

# --- New row 15 ---
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = "Vega Modelo de Temuco"
$ws.Range("C15").Value = "La Araucanía"
$ws.Range("D15").Value = [datetime]"2021-10-28"
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = 300000000
$ws.Range("G15").Value = "Espárragos"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Extra"
$ws.Range("J15").Value = 40
$ws.Range("K15").Value = 2000
$ws.Range("L15").Value = 2000
$ws.Range("M15").Value = 2000
$ws.Range("N15").Value = "`$/kilo"
$ws.Range("O15").Value = "Región del Maule"
$ws.Range("P15").Value = 2000
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = "Hortaliza"

# --- New row 16 ---
$ws.Range("A16").Value = 10
$ws.Range("B16").Value = "Vega Modelo de Temuco"
$ws.Range("C16").Value = "La Araucanía"
$ws.Range("D16").Value = [datetime]"2021-10-28"
$ws.Range("E16").Value = 9
$ws.Range("F16").Value = 300000000
$ws.Range("G16").Value = "Espárragos"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 550
$ws.Range("K16").Value = 1200
$ws.Range("L16").Value = 1300
$ws.Range("M16").Value = 1245
$ws.Range("N16").Value = "`$/kilo"
$ws.Range("O16").Value = "Región del Maule"
$ws.Range("P16").Value = 1245
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = "Hortaliza"
